$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.465.46"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "2.106.39"
$ws.Range("E3").Value = "  +1.47%  "
$ws.Range("E4").Value = "  -0.58%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "336.17"
$ws.Range("E5").Value = "  +2.39%  "
$ws.Range("E6").Value = "  -0.46%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5244"
$ws.Range("E7").Value = "  +0.87%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4587"
$ws.Range("E8").Value = "  +6.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.44"
$ws.Range("E9").Value = "  +15.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08948"
$ws.Range("E10").Value = "  +2.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.179"
$ws.Range("E11").Value = "  +2.49%  "
$ws.Range("E12").Value = "  +1.29%  "
$ws.Range("D13").Value = "2.081.02"
$ws.Range("E13").Value = "  -0.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.800"
$ws.Range("E14").Value = "  +2.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.964"
$ws.Range("E15").Value = "  +3.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "96.57"
$ws.Range("E16").Value = "  +1.17%  "
$ws.Range("E17").Value = "  -0.48%  "
$ws.Range("E18").Value = "  +1.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06633"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.31"
$ws.Range("E20").Value = "  +3.10%  "
$ws.Range("E21").Value = "  -0.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.305"
$ws.Range("E22").Value = "  +1.08%  "
$ws.Range("D23").Value = "30.518.53"
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("E24").Value = "  +0.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.366"
$ws.Range("E25").Value = "  +3.23%  "
$ws.Range("D26").Value = "2.336.03"
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("E27").Value = "  +1.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.568"
$ws.Range("E28").Value = "  +2.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "163.41"
$ws.Range("E29").Value = "  +0.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.80"
$ws.Range("E30").Value = "  +1.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.225"
$ws.Range("E31").Value = "  +4.14%  "
$ws.Range("B32").Value = "ARBITRUM"
$ws.Range("C32").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.707"
$ws.Range("E32").Value = "  +12.24%  "
$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.1074"
$ws.Range("E33").Value = "  +0.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.206"
$ws.Range("E34").Value = "  +2.94%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.926"
$ws.Range("E35").Value = "  +2.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.47"
$ws.Range("E36").Value = "  +8.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02578"
$ws.Range("E37").Value = "  +0.90%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06847"
$ws.Range("E38").Value = "  +3.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.558"
$ws.Range("E39").Value = "  +1.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.85"
$ws.Range("E40").Value = "  +3.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2305"
$ws.Range("E41").Value = "  +2.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6906"
$ws.Range("E42").Value = "  +3.09%  "
$ws.Range("E43").Value = "  +0.45%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.352"
$ws.Range("E44").Value = "  +7.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.000"
$ws.Range("E45").Value = "  -0.44%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "14.05"
$ws.Range("E46").Value = "  +1.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6388"
$ws.Range("E47").Value = "  +1.37%  "
$ws.Range("E48").Value = "  +1.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000351"
$ws.Range("E49").Value = "  +25.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.248"
$ws.Range("E50").Value = "  +1.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "83.63"
$ws.Range("E51").Value = "  +2.14%  "
